$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 612, shifting the existing data (old rows
# 612-660) down to rows 615-663. This also grows the sheet dimension from
# A1:R660 to A1:R663 automatically.
$ws.Rows("612:614").Insert()

# ---- New row 612 ----
$ws.Range("A612").Value2 = 10
$ws.Range("B612").Value2 = "Vega Modelo de Temuco"
$ws.Range("C612").Value2 = "La Araucanía"
$ws.Range("D612").Value2 = 44826
$ws.Range("D612").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E612").Value2 = 9
$ws.Range("F612").Value2 = 100112045
$ws.Range("G612").Value2 = "Zapallo"
$ws.Range("H612").Value2 = "Camote"
$ws.Range("I612").Value2 = "2a (guarda)"
$ws.Range("J612").Value2 = 300
$ws.Range("K612").Value2 = 1000
$ws.Range("L612").Value2 = 1000
$ws.Range("M612").Value2 = 1000
$ws.Range("N612").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O612").Value2 = "Región de O'Higgins"
$ws.Range("P612").Value2 = 1000
$ws.Range("Q612").Value2 = 1
$ws.Range("R612").Value2 = "Hortaliza"

# ---- New row 613 ----
$ws.Range("A613").Value2 = 10
$ws.Range("B613").Value2 = "Vega Modelo de Temuco"
$ws.Range("C613").Value2 = "La Araucanía"
$ws.Range("D613").Value2 = 44826
$ws.Range("D613").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E613").Value2 = 9
$ws.Range("F613").Value2 = 100112045
$ws.Range("G613").Value2 = "Zapallo"
$ws.Range("H613").Value2 = "Paine"
$ws.Range("I613").Value2 = "1a (guarda)"
$ws.Range("J613").Value2 = 800
$ws.Range("K613").Value2 = 500
$ws.Range("L613").Value2 = 500
$ws.Range("M613").Value2 = 500
$ws.Range("N613").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O613").Value2 = "Región de O'Higgins"
$ws.Range("P613").Value2 = 500
$ws.Range("Q613").Value2 = 1
$ws.Range("R613").Value2 = "Hortaliza"

# ---- New row 614 ----
$ws.Range("A614").Value2 = 10
$ws.Range("B614").Value2 = "Vega Modelo de Temuco"
$ws.Range("C614").Value2 = "La Araucanía"
$ws.Range("D614").Value2 = 44826
$ws.Range("D614").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E614").Value2 = 9
$ws.Range("F614").Value2 = 100112045
$ws.Range("G614").Value2 = "Zapallo"
$ws.Range("H614").Value2 = "Paine"
$ws.Range("I614").Value2 = "2a (guarda)"
$ws.Range("J614").Value2 = 300
$ws.Range("K614").Value2 = 400
$ws.Range("L614").Value2 = 400
$ws.Range("M614").Value2 = 400
$ws.Range("N614").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O614").Value2 = "Región de O'Higgins"
$ws.Range("P614").Value2 = 400
$ws.Range("Q614").Value2 = 1
$ws.Range("R614").Value2 = "Hortaliza"
